# Update column F ("dSF") values on Sheet1 to reflect repulled/recalculated data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -3
    4  = 1
    5  = -2
    6  = 4
    7  = -1
    8  = 2
    9  = -2
    10 = -1
    12 = -2
    13 = 2
    14 = 1
    15 = -2
    16 = 1
    17 = -1
    18 = 3
    19 = 1
    20 = -4
    21 = 2
    22 = 3
    23 = 3
    26 = 1
    27 = -6
    29 = 4
    30 = -7
    31 = -2
    32 = 6
    33 = -3
    35 = 2
    36 = 1
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 6).Value = $updates[$row]
}
